$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2537.6836
$ws.Range("J17").Value = 2537.6836
$ws.Range("L17").Value = 7613.050799999999
$ws.Range("N17").Value = -7949.050799999999

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2423
$ws.Range("I98").Value = 2201.25
$ws.Range("K98").Value = 2201.25
$ws.Range("M98").Value = -703.25

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 421.79413
$ws.Range("I103").Value = 455.75
$ws.Range("J103").Value = 417.26666
$ws.Range("K103").Value = 1367.25
$ws.Range("L103").Value = 1251.79998
$ws.Range("M103").Value = -781.25
$ws.Range("N103").Value = -2423.79998

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 22728124
$ws.Range("I107").Value = 23810274
$ws.Range("J107").Value = 2999
$ws.Range("K107").Value = 23810274
$ws.Range("L107").Value = 2999
$ws.Range("M107").Value = -23808354
$ws.Range("N107").Value = -6839

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2685.9412
$ws.Range("J112").Value = 2685.9412
$ws.Range("L112").Value = 8057.823600000001
$ws.Range("N112").Value = -10273.8236

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2423
$ws.Range("I122").Value = 2201.25
$ws.Range("K122").Value = 6603.75
$ws.Range("M122").Value = -4153.75

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 9240.796
$ws.Range("I137").Value = 3387.92
$ws.Range("J137").Value = 15337.542
$ws.Range("K137").Value = 10163.76
$ws.Range("L137").Value = 46012.626
$ws.Range("M137").Value = -7613.76
$ws.Range("N137").Value = -51112.626

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3545.7563
$ws.Range("I138").Value = 3224.2646
$ws.Range("J138").Value = 3794.182
$ws.Range("K138").Value = 9672.793799999999
$ws.Range("L138").Value = 11382.546
$ws.Range("M138").Value = -4532.793799999999
$ws.Range("N138").Value = -21662.546

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3146084.8
$ws.Range("I32").Value = 3300023.2
$ws.Range("K32").Value = 3300023.2
$ws.Range("M32").Value = -3299736.2

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1920.8125
$ws.Range("I45").Value = 1711
$ws.Range("J45").Value = 1969.2307
$ws.Range("K45").Value = 1711
$ws.Range("L45").Value = 1969.2307
$ws.Range("M45").Value = -1334
$ws.Range("N45").Value = -2723.2307

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 14708641
$ws.Range("I110").Value = 22729354
$ws.Range("K110").Value = 22729354
$ws.Range("M110").Value = -22727309

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1945
$ws.Range("I122").Value = 1822.75
$ws.Range("J122").Value = 2042.8
$ws.Range("K122").Value = 5468.25
$ws.Range("L122").Value = 6128.4
$ws.Range("M122").Value = -3018.25
$ws.Range("N122").Value = -11028.4

# BSM row 26
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 173500
$ws.Range("I26").Value = 173500
$ws.Range("K26").Value = 173500
$ws.Range("M26").Value = -173208

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26320426
$ws.Range("I31").Value = 83335090
$ws.Range("J31").Value = 5965.769
$ws.Range("K31").Value = 83335090
$ws.Range("L31").Value = 5965.769
$ws.Range("M31").Value = -83334795
$ws.Range("N31").Value = -6555.769

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 26320426
$ws.Range("I34").Value = 83335090
$ws.Range("J34").Value = 5965.769
$ws.Range("K34").Value = 83335090
$ws.Range("L34").Value = 5965.769
$ws.Range("M34").Value = -83334888
$ws.Range("N34").Value = -6369.769

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 50148.555
$ws.Range("I41").Value = 47725.668
$ws.Range("J41").Value = 51360
$ws.Range("K41").Value = 47725.668
$ws.Range("L41").Value = 51360
$ws.Range("M41").Value = -47297.668
$ws.Range("N41").Value = -52216

# CRP row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 39141
$ws.Range("J50").Value = 39141
$ws.Range("L50").Value = 39141
$ws.Range("N50").Value = -40391

# CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 39047.332
$ws.Range("I51").Value = 2090
$ws.Range("J51").Value = 46438.8
$ws.Range("K51").Value = 2090
$ws.Range("L51").Value = 46438.8
$ws.Range("M51").Value = -1354
$ws.Range("N51").Value = -47910.8

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4262.08
$ws.Range("I58").Value = 4005.4194
$ws.Range("J58").Value = 4680.8423
$ws.Range("K58").Value = 4005.4194
$ws.Range("L58").Value = 4680.8423
$ws.Range("M58").Value = -3802.4194
$ws.Range("N58").Value = -5086.8423

# CRP row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 29966.666
$ws.Range("J60").Value = 29966.666
$ws.Range("L60").Value = 29966.666
$ws.Range("N60").Value = -30988.666

# CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 39047.332
$ws.Range("I61").Value = 2090
$ws.Range("J61").Value = 46438.8
$ws.Range("K61").Value = 2090
$ws.Range("L61").Value = 46438.8
$ws.Range("M61").Value = -1742
$ws.Range("N61").Value = -47134.8

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9152.462
$ws.Range("J99").Value = 9998.666999999999
$ws.Range("L99").Value = 9998.666999999999
$ws.Range("N99").Value = -12994.667

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 9152.462
$ws.Range("J126").Value = 9998.666999999999
$ws.Range("L126").Value = 29996.001
$ws.Range("N126").Value = -34936.001

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2819.0576
$ws.Range("J132").Value = 4291.8335
$ws.Range("L132").Value = 12875.5005
$ws.Range("N132").Value = -17935.5005

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4262.08
$ws.Range("I136").Value = 4005.4194
$ws.Range("J136").Value = 4680.8423
$ws.Range("K136").Value = 12016.2582
$ws.Range("L136").Value = 14042.5269
$ws.Range("M136").Value = -9466.2582
$ws.Range("N136").Value = -19142.5269

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2146022
$ws.Range("I4").Value = 781967.75
$ws.Range("K4").Value = 2345903.25
$ws.Range("M4").Value = -2345791.25

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1500
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# CUL row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1666.3334
$ws.Range("J97").Value = 1500
$ws.Range("L97").Value = 4500
$ws.Range("N97").Value = -5492

# CUL row 100
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 2725
$ws.Range("I100").Value = 966.6667
$ws.Range("J100").Value = 8000
$ws.Range("K100").Value = 2900.0001
$ws.Range("L100").Value = 24000
$ws.Range("M100").Value = -2089.0001
$ws.Range("N100").Value = -25622

# CUL row 109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 828.5
$ws.Range("I109").Value = 907.2222
$ws.Range("K109").Value = 2721.6666
$ws.Range("M109").Value = -1681.6666

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1810.8462
$ws.Range("J113").Value = 1868.6842
$ws.Range("L113").Value = 5606.0526
$ws.Range("N113").Value = -9946.052599999999

# CUL row 124
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 3576.3333
$ws.Range("I124").Value = 3576.3333
$ws.Range("K124").Value = 10728.9999
$ws.Range("M124").Value = -5818.999899999999

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3548.5
$ws.Range("I132").Value = 4387
$ws.Range("J132").Value = 1368.4
$ws.Range("K132").Value = 39483
$ws.Range("L132").Value = 12315.6
$ws.Range("M132").Value = -36953
$ws.Range("N132").Value = -17375.6

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 50330.15
$ws.Range("I2").Value = 465.16666
$ws.Range("K2").Value = 465.16666
$ws.Range("M2").Value = -352.16666

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3292.158
$ws.Range("J132").Value = 7246.923
$ws.Range("L132").Value = 21740.769
$ws.Range("N132").Value = -26800.769

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 173384.5
$ws.Range("I7").Value = 173384.5
$ws.Range("K7").Value = 173384.5
$ws.Range("M7").Value = -173272.5

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 55559668
$ws.Range("J122").Value = 4498.5
$ws.Range("L122").Value = 13495.5
$ws.Range("N122").Value = -18395.5

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 173384.5
$ws.Range("I126").Value = 173384.5
$ws.Range("K126").Value = 520153.5
$ws.Range("M126").Value = -517683.5

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2302.7576
$ws.Range("I132").Value = 1676.3667
$ws.Range("K132").Value = 5029.1001
$ws.Range("M132").Value = -2499.1001

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 765.3570999999999
$ws.Range("I100").Value = 492.91666
$ws.Range("J100").Value = 2400
$ws.Range("K100").Value = 985.83332
$ws.Range("L100").Value = 4800
$ws.Range("M100").Value = -444.83332
$ws.Range("N100").Value = -5882

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1452.875
$ws.Range("I107").Value = 814.8570999999999
$ws.Range("K107").Value = 2444.5713
$ws.Range("M107").Value = -524.5712999999996

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6000
$ws.Range("I122").Value = 6000
$ws.Range("K122").Value = 18000
$ws.Range("M122").Value = -15550

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5796.154
$ws.Range("I132").Value = 3872.7222
$ws.Range("K132").Value = 11618.1666
$ws.Range("M132").Value = -9088.1666
